$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, and whether the cell must be
# forced to Text format first so Excel does not reinterpret a
# numeric-looking string (e.g. "1.008") as a real number.
$updates = @(
    @{ Cell = "D2"; Value = '26.383.45'; AsText = $false }
    @{ Cell = "E2"; Value = '  -2.40%  '; AsText = $false }
    @{ Cell = "D3"; Value = '1.794.17'; AsText = $false }
    @{ Cell = "E3"; Value = '  -2.55%  '; AsText = $false }
    @{ Cell = "D4"; Value = '1.008'; AsText = $true }
    @{ Cell = "E4"; Value = '  +0.19%  '; AsText = $false }
    @{ Cell = "E5"; Value = '  +0.15%  '; AsText = $false }
    @{ Cell = "D6"; Value = '306.92'; AsText = $true }
    @{ Cell = "E6"; Value = '  -2.05%  '; AsText = $false }
    @{ Cell = "D7"; Value = '0.4547'; AsText = $true }
    @{ Cell = "E7"; Value = '  -1.79%  '; AsText = $false }
    @{ Cell = "D8"; Value = '0.3626'; AsText = $true }
    @{ Cell = "E8"; Value = '  -1.60%  '; AsText = $false }
    @{ Cell = "B9"; Value = 'Dogecoin'; AsText = $false }
    @{ Cell = "C9"; Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; AsText = $false }
    @{ Cell = "D9"; Value = '0.07066'; AsText = $true }
    @{ Cell = "E9"; Value = '  -2.50%  '; AsText = $false }
    @{ Cell = "B10"; Value = 'Polygon'; AsText = $false }
    @{ Cell = "C10"; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; AsText = $false }
    @{ Cell = "D10"; Value = '0.8713'; AsText = $true }
    @{ Cell = "E10"; Value = '  -1.32%  '; AsText = $false }
    @{ Cell = "B11"; Value = 'TRON'; AsText = $false }
    @{ Cell = "C11"; Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; AsText = $false }
    @{ Cell = "D11"; Value = '0.07794'; AsText = $true }
    @{ Cell = "E11"; Value = '  -0.51%  '; AsText = $false }
    @{ Cell = "B12"; Value = 'Solana'; AsText = $false }
    @{ Cell = "C12"; Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; AsText = $false }
    @{ Cell = "D12"; Value = '19.38'; AsText = $true }
    @{ Cell = "E12"; Value = '  -1.92%  '; AsText = $false }
    @{ Cell = "B13"; Value = 'WrappedEther'; AsText = $false }
    @{ Cell = "C13"; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; AsText = $false }
    @{ Cell = "D13"; Value = '1.753.99'; AsText = $false }
    @{ Cell = "E13"; Value = '  -5.36%  '; AsText = $false }
    @{ Cell = "B14"; Value = 'Polkadot'; AsText = $false }
    @{ Cell = "C14"; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; AsText = $false }
    @{ Cell = "D14"; Value = '5.260'; AsText = $true }
    @{ Cell = "E14"; Value = '  -2.01%  '; AsText = $false }
    @{ Cell = "B15"; Value = 'Chainlink'; AsText = $false }
    @{ Cell = "C15"; Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; AsText = $false }
    @{ Cell = "D15"; Value = '6.311'; AsText = $true }
    @{ Cell = "E15"; Value = '  -2.49%  '; AsText = $false }
    @{ Cell = "B16"; Value = 'Litecoin'; AsText = $false }
    @{ Cell = "C16"; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; AsText = $false }
    @{ Cell = "D16"; Value = '84.61'; AsText = $true }
    @{ Cell = "E16"; Value = '  -6.87%  '; AsText = $false }
    @{ Cell = "B17"; Value = 'BinanceUSD'; AsText = $false }
    @{ Cell = "C17"; Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; AsText = $false }
    @{ Cell = "D17"; Value = '1.009'; AsText = $true }
    @{ Cell = "E17"; Value = '  +0.26%  '; AsText = $false }
    @{ Cell = "B18"; Value = 'ShibaInu'; AsText = $false }
    @{ Cell = "C18"; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; AsText = $false }
    @{ Cell = "D18"; Value = '0.000008508'; AsText = $true }
    @{ Cell = "E18"; Value = '  -3.46%  '; AsText = $false }
    @{ Cell = "B19"; Value = 'Dai'; AsText = $false }
    @{ Cell = "C19"; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; AsText = $false }
    @{ Cell = "D19"; Value = '1.008'; AsText = $true }
    @{ Cell = "E19"; Value = '  +0.20%  '; AsText = $false }
    @{ Cell = "B20"; Value = 'WrappedBTC'; AsText = $false }
    @{ Cell = "C20"; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; AsText = $false }
    @{ Cell = "D20"; Value = '26.440.28'; AsText = $false }
    @{ Cell = "E20"; Value = '  -2.35%  '; AsText = $false }
    @{ Cell = "B21"; Value = 'Avalanche'; AsText = $false }
    @{ Cell = "C21"; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; AsText = $false }
    @{ Cell = "D21"; Value = '14.15'; AsText = $true }
    @{ Cell = "E21"; Value = '  -3.03%  '; AsText = $false }
    @{ Cell = "B22"; Value = 'Uniswap'; AsText = $false }
    @{ Cell = "C22"; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; AsText = $false }
    @{ Cell = "D22"; Value = '4.971'; AsText = $true }
    @{ Cell = "E22"; Value = '  -1.05%  '; AsText = $false }
    @{ Cell = "B23"; Value = 'Cosmos'; AsText = $false }
    @{ Cell = "C23"; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; AsText = $false }
    @{ Cell = "D23"; Value = '10.48'; AsText = $true }
    @{ Cell = "E23"; Value = '  -0.17%  '; AsText = $false }
    @{ Cell = "D24"; Value = '1.997.47'; AsText = $false }
    @{ Cell = "E24"; Value = '  -6.64%  '; AsText = $false }
    @{ Cell = "B25"; Value = 'Toncoin'; AsText = $false }
    @{ Cell = "C25"; Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; AsText = $false }
    @{ Cell = "D25"; Value = '1.979'; AsText = $true }
    @{ Cell = "E25"; Value = '  -1.28%  '; AsText = $false }
    @{ Cell = "B26"; Value = 'Monero'; AsText = $false }
    @{ Cell = "C26"; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; AsText = $false }
    @{ Cell = "D26"; Value = '152.04'; AsText = $true }
    @{ Cell = "E26"; Value = '  +0.81%  '; AsText = $false }
    @{ Cell = "B27"; Value = 'EthereumClassic'; AsText = $false }
    @{ Cell = "C27"; Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; AsText = $false }
    @{ Cell = "D27"; Value = '17.83'; AsText = $true }
    @{ Cell = "E27"; Value = '  -2.85%  '; AsText = $false }
    @{ Cell = "B28"; Value = 'LidoDAOToken'; AsText = $false }
    @{ Cell = "C28"; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; AsText = $false }
    @{ Cell = "D28"; Value = '2.034'; AsText = $true }
    @{ Cell = "E28"; Value = '  +1.00%  '; AsText = $false }
    @{ Cell = "B29"; Value = 'BitcoinCash'; AsText = $false }
    @{ Cell = "C29"; Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; AsText = $false }
    @{ Cell = "D29"; Value = '112.05'; AsText = $true }
    @{ Cell = "E29"; Value = '  -2.91%  '; AsText = $false }
    @{ Cell = "B30"; Value = 'InternetComputer(DFINITY)'; AsText = $false }
    @{ Cell = "C30"; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; AsText = $false }
    @{ Cell = "D30"; Value = '4.824'; AsText = $true }
    @{ Cell = "E30"; Value = '  -2.98%  '; AsText = $false }
    @{ Cell = "B31"; Value = 'Stellar'; AsText = $false }
    @{ Cell = "C31"; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; AsText = $false }
    @{ Cell = "D31"; Value = '0.08659'; AsText = $true }
    @{ Cell = "E31"; Value = '  -2.21%  '; AsText = $false }
    @{ Cell = "B32"; Value = 'HuobiToken'; AsText = $false }
    @{ Cell = "C32"; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; AsText = $false }
    @{ Cell = "D32"; Value = '3.025'; AsText = $true }
    @{ Cell = "E32"; Value = '  -3.80%  '; AsText = $false }
    @{ Cell = "B33"; Value = 'Filecoin'; AsText = $false }
    @{ Cell = "C33"; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; AsText = $false }
    @{ Cell = "D33"; Value = '4.435'; AsText = $true }
    @{ Cell = "E33"; Value = '  -1.13%  '; AsText = $false }
    @{ Cell = "B34"; Value = 'ImmutableX'; AsText = $false }
    @{ Cell = "C34"; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; AsText = $false }
    @{ Cell = "D34"; Value = '0.7152'; AsText = $true }
    @{ Cell = "E34"; Value = '  -6.88%  '; AsText = $false }
    @{ Cell = "B35"; Value = 'RenderToken'; AsText = $false }
    @{ Cell = "C35"; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; AsText = $false }
    @{ Cell = "D35"; Value = '2.668'; AsText = $true }
    @{ Cell = "E35"; Value = '  +0.19%  '; AsText = $false }
    @{ Cell = "D36"; Value = '1.110'; AsText = $true }
    @{ Cell = "E36"; Value = '  -2.74%  '; AsText = $false }
    @{ Cell = "B37"; Value = 'Frax'; AsText = $false }
    @{ Cell = "C37"; Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'; AsText = $false }
    @{ Cell = "D37"; Value = '1.006'; AsText = $true }
    @{ Cell = "E37"; Value = '  +0.09%  '; AsText = $false }
    @{ Cell = "B38"; Value = 'TrustWalletToken'; AsText = $false }
    @{ Cell = "C38"; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; AsText = $false }
    @{ Cell = "D38"; Value = '1.079'; AsText = $true }
    @{ Cell = "E38"; Value = '  -1.78%  '; AsText = $false }
    @{ Cell = "B39"; Value = 'VeChain'; AsText = $false }
    @{ Cell = "C39"; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; AsText = $false }
    @{ Cell = "D39"; Value = '0.01938'; AsText = $true }
    @{ Cell = "E39"; Value = '  +0.17%  '; AsText = $false }
    @{ Cell = "B40"; Value = 'Hedera'; AsText = $false }
    @{ Cell = "C40"; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; AsText = $false }
    @{ Cell = "D40"; Value = '0.05082'; AsText = $true }
    @{ Cell = "E40"; Value = '  -1.99%  '; AsText = $false }
    @{ Cell = "B41"; Value = 'MXToken'; AsText = $false }
    @{ Cell = "C41"; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; AsText = $false }
    @{ Cell = "D41"; Value = '2.863'; AsText = $true }
    @{ Cell = "E41"; Value = '  -2.37%  '; AsText = $false }
    @{ Cell = "B42"; Value = 'FraxShare'; AsText = $false }
    @{ Cell = "C42"; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; AsText = $false }
    @{ Cell = "D42"; Value = '6.888'; AsText = $true }
    @{ Cell = "E42"; Value = '  -1.48%  '; AsText = $false }
    @{ Cell = "B43"; Value = 'TheSandbox'; AsText = $false }
    @{ Cell = "C43"; Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; AsText = $false }
    @{ Cell = "D43"; Value = '0.4983'; AsText = $true }
    @{ Cell = "E43"; Value = '  -0.86%  '; AsText = $false }
    @{ Cell = "B44"; Value = 'Algorand'; AsText = $false }
    @{ Cell = "C44"; Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; AsText = $false }
    @{ Cell = "D44"; Value = '0.1515'; AsText = $true }
    @{ Cell = "E44"; Value = '  -5.75%  '; AsText = $false }
    @{ Cell = "B45"; Value = 'Aptos'; AsText = $false }
    @{ Cell = "C45"; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; AsText = $false }
    @{ Cell = "D45"; Value = '7.964'; AsText = $true }
    @{ Cell = "E45"; Value = '  -5.56%  '; AsText = $false }
    @{ Cell = "B46"; Value = 'PaxDollar'; AsText = $false }
    @{ Cell = "C46"; Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; AsText = $false }
    @{ Cell = "D46"; Value = '1.009'; AsText = $true }
    @{ Cell = "E46"; Value = '  +0.24%  '; AsText = $false }
    @{ Cell = "B47"; Value = 'Decentraland'; AsText = $false }
    @{ Cell = "C47"; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; AsText = $false }
    @{ Cell = "D47"; Value = '0.4583'; AsText = $true }
    @{ Cell = "E47"; Value = '  -2.84%  '; AsText = $false }
    @{ Cell = "B48"; Value = 'EnergySwap'; AsText = $false }
    @{ Cell = "C48"; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; AsText = $false }
    @{ Cell = "D48"; Value = '9.859'; AsText = $true }
    @{ Cell = "E48"; Value = '  -3.75%  '; AsText = $false }
    @{ Cell = "B49"; Value = 'Quant'; AsText = $false }
    @{ Cell = "C49"; Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; AsText = $false }
    @{ Cell = "D49"; Value = '99.48'; AsText = $true }
    @{ Cell = "E49"; Value = '  -3.25%  '; AsText = $false }
    @{ Cell = "B50"; Value = 'NEARProtocol'; AsText = $false }
    @{ Cell = "C50"; Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; AsText = $false }
    @{ Cell = "D50"; Value = '1.581'; AsText = $true }
    @{ Cell = "E50"; Value = '  -2.38%  '; AsText = $false }
    @{ Cell = "B51"; Value = 'Cronos'; AsText = $false }
    @{ Cell = "C51"; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; AsText = $false }
    @{ Cell = "D51"; Value = '0.05966'; AsText = $true }
    @{ Cell = "E51"; Value = '  -2.84%  '; AsText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.AsText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
